$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Testing"
$ws.Range("C9").Value = "Testing"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
